$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 242, shifting existing rows 242-359 down to 243-360
$ws.Rows.Item(242).Insert()

# Populate the newly inserted row 242 with the new record's data
$ws.Cells.Item(242, 1).Value = 10
$ws.Cells.Item(242, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(242, 3).Value = "La Araucanía"
$ws.Cells.Item(242, 4).Value = "2022-02-17"
$ws.Cells.Item(242, 5).Value = 9
$ws.Cells.Item(242, 6).Value = 100112032
$ws.Cells.Item(242, 7).Value = "Zapallo italiano"
$ws.Cells.Item(242, 8).Value = "Sin especificar"
$ws.Cells.Item(242, 9).Value = "Primera"
$ws.Cells.Item(242, 10).Value = 400
$ws.Cells.Item(242, 11).Value = 10000
$ws.Cells.Item(242, 12).Value = 12000
$ws.Cells.Item(242, 13).Value = 11000
$ws.Cells.Item(242, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(242, 15).Value = "Región del Maule"
$ws.Cells.Item(242, 16).Value = 183
$ws.Cells.Item(242, 17).Value = 60
$ws.Cells.Item(242, 18).Value = "Hortaliza"
